Write-Host "no-op at all"
